$d = $word.ActiveDocument

# --- Locate the "Week 7" block of three list paragraphs by their text ---
$autoDrawIdx = 0
$continueIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Automatically draw a new Event Card*") {
        $autoDrawIdx = $i
    }
    if ($t -like "Continue implementing event cards*" -and $autoDrawIdx -ne 0 -and $continueIdx -eq 0) {
        $continueIdx = $i
    }
}

# 1) Add green highlighting to the "Automatically draw..." paragraph
#    (covers both the run and the paragraph mark, matching the source rPr).
$autoPara = $d.Paragraphs.Item($autoDrawIdx)
$autoPara.Range.Font.HighlightColorIndex = 4

# 2) Remove the "Continue implementing event cards (3 - 4 more)" paragraph
#    from its old position (right after "Automatically draw...").
$continuePara = $d.Paragraphs.Item($continueIdx)
$continuePara.Range.Delete()

# 3) The paragraph that used to hold "Test event cards..." (with the
#    _GoBack bookmark) is now at $continueIdx; replace its text in place so
#    the bookmark stays put, re-creating the original two runs.
$targetPara = $d.Paragraphs.Item($continueIdx)
$targetPara.Range.Text = "Continue implementing event cards"
$afterFirstRun = $d.Range($targetPara.Range.End - 1, $targetPara.Range.End - 1)
$afterFirstRun.InsertAfter(" (3 – 4 more)")

# 4) Insert a brand-new paragraph after it (no bookmark) carrying the
#    "Test event cards and resulting business logic" text that used to
#    live in the bookmarked paragraph.
$targetPara = $d.Paragraphs.Item($continueIdx)
$insertPoint = $d.Range($targetPara.Range.End - 1, $targetPara.Range.End - 1)
$insertPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($continueIdx + 1)
$newPara.Range.Text = "Test event cards and resulting business logic"
